$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 title cells: drop the period after "6.4.1.2" in the Russian / English titles.
$ws.Range("B1").Value = "6.4.1.2 Потери воды при транспортировке"
$ws.Range("C1").Value = "6.4.1.2 Percentage of water loss during transportation"

# Updated data values (column P, 2022).
$ws.Range("P5").Value = 2388
$ws.Range("P10").Value = 335.3
$ws.Range("P16").Value = 27.3
$ws.Range("P21").Value = 24.3

# Move the active selection from R9 to S3.
$ws.Range("S3").Select()
